# feat: add 2022-Q3 data
#
# The workbook has a "总计" (totals) sheet and a single quarterly sheet
# named "2022-Q2". This change:
#   1. Preserves the existing "2022-Q2" sheet (and its data) by copying it
#      to a new sheet placed right after it.
#   2. Replaces the data on the original sheet with the new 2022-Q3 figures
#      and renames it to "2022-Q3".
#   3. Renames the copy back to "2022-Q2" so history is kept intact.
#   4. Updates the "总计" summary sheet: the existing 2022-Q2 row becomes the
#      2022-Q3 row (with the new total market value) and a new row is
#      appended below it with the original 2022-Q2 figures.

$xlPasteValues = -4163
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

# Helper: write $value into $target as a genuine text cell (matching the
# source file's convention of storing figures such as "22.76" as text, not
# numbers) without picking up a one-off "@" text number-format/style. A
# scratch cell gets a literal-text formula; copying its computed value with
# "paste values" carries over the string *type* with no formatting at all.
function Set-TextValue($targetSheet, $targetAddr, $scratchAddr, $value) {
    $targetSheet.Range($scratchAddr).Formula = '="' + $value + '"'
    $targetSheet.Range($scratchAddr).Copy()
    $targetSheet.Range($targetAddr).PasteSpecial($xlPasteValues)
    $targetSheet.Range($scratchAddr).Clear()
}

# --- 1. Duplicate the current "2022-Q2" sheet so its data survives -------
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item(3)

# --- 2. Clear the original sheet and fill in the new 2022-Q3 figures -----
$q2.Cells.Clear()

$q2.Range("B1").Value = "基金代码"
$q2.Range("C1").Value = "基金名称"
$q2.Range("D1").Value = "基金规模"
$q2.Range("E1").Value = "股票总仓位"
$q2.Range("F1").Value = "仓位占比"
$q2.Range("G1").Value = "持有市值(亿元)"
$q2.Range("H1").Value = "仓位排名"

$q2.Range("A2").Value = 0
Set-TextValue $q2 "B2" "Z1" "377016"
$q2.Range("C2").Value = "上投摩根亚太优势混合（QDII）"
Set-TextValue $q2 "D2" "Z1" "22.76"
Set-TextValue $q2 "E2" "Z1" "90.69"
Set-TextValue $q2 "F2" "Z1" "1.95"
Set-TextValue $q2 "G2" "Z1" "0.4438"
$q2.Range("H2").Value = 10

# Match the bold/bordered/centered header style used elsewhere in the
# workbook (e.g. the "总计" sheet header).
$total.Range("B1").Copy()
$q2.Range("B1:H1").PasteSpecial($xlPasteFormats)
$total.Range("A2").Copy()
$q2.Range("A2").PasteSpecial($xlPasteFormats)

$q2.PageSetup.LeftMargin = 54
$q2.PageSetup.RightMargin = 54
$q2.PageSetup.TopMargin = 72
$q2.PageSetup.BottomMargin = 72
$q2.PageSetup.HeaderMargin = 36
$q2.PageSetup.FooterMargin = 36

# --- 3. Rename sheets (rename the source away from "2022-Q2" first so the
#        copy can reclaim that name without a collision) -----------------
$q2.Name = "2022-Q3"
$q2Copy.Name = "2022-Q2"

# --- 4. Update the "总计" summary sheet -----------------------------------
$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.44

$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial($xlPasteFormats)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.55
